# "ignoring decommissioned power plants"
# Update the "Import Priorities" sheet: drop the now-unused
# ConventionalPlantOperator / VariableRenewableOperator / electricity rows,
# bump "unit" priority, and add PowerPlantsInstalled / Decommissioned rows
# at the bottom (lowest priority), replacing TechnologyPotentials' old slot.

$wb = $excel.ActiveWorkbook

$wsImport = $wb.Worksheets.Item("Import Priorities")
$wsCoupling = $wb.Worksheets.Item("Coupling Parameters")

# Remove the three rows that are no longer relevant for import priority
# (they currently sit at rows 7-9: ConventionalPlantOperator, VariableRenewableOperator, electricity).
$wsImport.Rows.Item(7).Delete()
$wsImport.Rows.Item(7).Delete()
$wsImport.Rows.Item(7).Delete()

# "unit" now moves up to row 7 - bump its priority from 3 to 4.
$wsImport.Range("B7").Value = 4

# "TechnologyPotentials" now sits at row 8 - bump its priority from 1 to 3.
$wsImport.Range("B8").Value = 3

# Append the two new, lower-priority rows. Write the shared strings in the
# same order they appear in the target workbook (PowerPlantsInstalled then
# Decommissioned) so the new <si> entries land in the expected order.
$wsImport.Range("A10").Value = "PowerPlantsInstalled"
$wsImport.Range("B10").Value = 1

$wsImport.Range("A9").Value = "Decommissioned"
$wsImport.Range("B9").Value = 2

# Make "Import Priorities" the active sheet/tab with the selection left at G7,
# and leave "Coupling Parameters" selected at D15 (no longer the active tab).
$wsCoupling.Range("D15").Select() | Out-Null
$wsImport.Activate()
$wsImport.Range("G7").Select() | Out-Null
